$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.270.31"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "'1.928.23"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'248.68"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "'0.7163"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'27.76"
$ws.Range("E8").Value = "  -1.88%  "
$ws.Range("D9").Value = "'0.3184"
$ws.Range("E9").Value = "  -4.62%  "
$ws.Range("D10").Value = "'0.07088"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").Value = "'0.7907"
$ws.Range("E11").Value = "  -2.92%  "
$ws.Range("D12").Value = "'0.07985"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").Value = "'1.922.61"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "'5.383"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").Value = "'94.86"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "'14.66"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").Value = "'30.265.23"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "'256.39"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").Value = "'0.000008047"
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").Value = "'5.759"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").Value = "'2.180.08"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "'1.0000"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'6.826"
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("D25").Value = "'9.529"
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("D26").Value = "'165.31"
$ws.Range("E26").Value = "  +3.39%  "
$ws.Range("D27").Value = "'19.08"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").Value = "'2.258"
$ws.Range("E28").Value = "  -6.97%  "
$ws.Range("D29").Value = "'0.1260"
$ws.Range("E29").Value = "  -5.90%  "
$ws.Range("D30").Value = "'1.352"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").Value = "'1.525"
$ws.Range("E31").Value = "  -2.11%  "
$ws.Range("D32").Value = "'4.394"
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").Value = "'4.121"
$ws.Range("E33").Value = "  -2.39%  "
$ws.Range("D34").Value = "'0.05134"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").Value = "'1.273"
$ws.Range("E35").Value = "  +3.05%  "
$ws.Range("D36").Value = "'0.7455"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "'2.763"
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("D38").Value = "'0.01961"
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("D39").Value = "'2.799"
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("D40").Value = "'78.30"
$ws.Range("E40").Value = "  -1.00%  "
$ws.Range("D41").Value = "'6.363"
$ws.Range("E41").Value = "  -4.46%  "
$ws.Range("D42").Value = "'0.4502"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").Value = "'1.985"
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("D44").Value = "'0.8465"
$ws.Range("E44").Value = "  +0.92%  "
$ws.Range("D45").Value = "'0.9996"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "'100.47"
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("D47").Value = "'9.754"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").Value = "'7.412"
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("D49").Value = "'36.63"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("D50").Value = "'0.06109"
$ws.Range("E50").Value = "  +2.39%  "
$ws.Range("D51").Value = "'0.4186"
$ws.Range("E51").Value = "  +1.96%  "
